# Update the "Notes" sheet:
#  - Row 2 (Description): reword.
#  - Row 4 (Source): reword, and insert a new "Source-link" row right after it.
#  - The old "as-is / open-use license" line (originally row 13, now shifted
#    down by the inserted row above) is reworded, and a new "More information
#    on licensing..." row is inserted right after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# --- Description / Source updates -----------------------------------------
$ws.Cells.Item(2, 1).Value = "Description: Population Density (Pop. Per Sq Km)"
$ws.Cells.Item(4, 1).Value = "Source: National population and Housing census 2014: Provisional Results - Uganda Bureau of Statistics"

# Insert a new row right after the Source row (row 4) for the source link.
$ws.Rows(5).Insert()
$ws.Cells.Item(5, 1).Value = "Source-link: http://www.ubos.org/onlinefiles/uploads/ubos/NPHC/NPHC%202014%20PROVISIONAL%20RESULTS%20REPORT.pdf"

# --- License updates --------------------------------------------------------
# After the insert above, the old row 13 ("It is provided on an as-is basis
# under an open-use license.") is now row 14.
$ws.Cells.Item(14, 1).Value = "It is licensed under a Creative Commons Attribution 4.0 International license."

# Insert a new row right after it for the licensing info link.
$ws.Rows(15).Insert()
$ws.Cells.Item(15, 1).Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"
